$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '56.496.94'
Set-TextValue "E2" '  -4.12%  '
Set-TextValue "D3" '2.376.46'
Set-TextValue "E3" '  -5.10%  '
Set-TextValue "E4" '  +0.02%  '
Set-TextValue "D5" '502.57'
Set-TextValue "E5" '  -6.09%  '
Set-TextValue "D6" '128.96'
Set-TextValue "E6" '  -3.98%  '
Set-TextValue "D7" '0.997'
Set-TextValue "E7" '  -0.32%  '
Set-TextValue "E8" '  -2.88%  '
Set-TextValue "D9" '2.398.68'
Set-TextValue "E9" '  -4.34%  '
Set-TextValue "D10" '0.0957'
Set-TextValue "E10" '  -3.93%  '
Set-TextValue "E11" '  -1.56%  '
Set-TextValue "D12" '0.318'
Set-TextValue "E12" '  -3.57%  '
Set-TextValue "D13" '4.61'
Set-TextValue "E13" '  -10.76%  '
Set-TextValue "D14" '2.802.17'
Set-TextValue "E14" '  -4.94%  '
Set-TextValue "D15" '56.958.51'
Set-TextValue "E15" '  -3.00%  '
Set-TextValue "D16" '21.58'
Set-TextValue "E16" '  -3.62%  '
Set-TextValue "E17" '  -3.38%  '
Set-TextValue "D18" '2.375.34'
Set-TextValue "E18" '  -5.25%  '
Set-TextValue "D19" '10.14'
Set-TextValue "E19" '  -4.67%  '
Set-TextValue "D20" '310.00'
Set-TextValue "E20" '  -3.47%  '
Set-TextValue "E21" '  -5.42%  '
Set-TextValue "D22" '6.21'
Set-TextValue "E22" '  -0.16%  '
Set-TextValue "D23" '0.998'
Set-TextValue "E23" '  -0.09%  '
Set-TextValue "D24" '65.43'
Set-TextValue "E24" '  -0.58%  '
Set-TextValue "D25" '1.00'
Set-TextValue "E25" '  +0.36%  '
Set-TextValue "D26" '2.493.17'
Set-TextValue "E26" '  -4.85%  '
Set-TextValue "D27" '0.372'
Set-TextValue "E27" '  -9.05%  '
Set-TextValue "E28" '  -6.00%  '
Set-TextValue "D29" '7.23'
Set-TextValue "E29" '  -3.03%  '
Set-TextValue "D30" '174.09'
Set-TextValue "E30" '  +1.27%  '
Set-TextValue "E31" '  -4.19%  '
Set-TextValue "D32" '0.0₃0711'
Set-TextValue "E32" '  -5.95%  '
Set-TextValue "D33" '6.11'
Set-TextValue "E33" '  -2.69%  '
Set-TextValue "D34" '0.998'
Set-TextValue "E34" '  -0.09%  '
Set-TextValue "E35" '  -7.32%  '
Set-TextValue "E36" '  -0.27%  '
Set-TextValue "E37" '  -1.90%  '
Set-TextValue "E38" '  -1.66%  '
Set-TextValue "D39" '3.75'
Set-TextValue "E39" '  -5.29%  '
Set-TextValue "D40" '35.85'
Set-TextValue "E40" '  -1.84%  '
Set-TextValue "E41" '  -6.06%  '
Set-TextValue "D42" '0.771'
Set-TextValue "E42" '  -7.15%  '
Set-TextValue "D43" '130.52'
Set-TextValue "E43" '  -0.50%  '
Set-TextValue "E44" '  -3.68%  '
Set-TextValue "D45" '4.87'
Set-TextValue "E45" '  -3.15%  '
Set-TextValue "D46" '0.573'
Set-TextValue "E46" '  -3.06%  '
Set-TextValue "D47" '253.98'
Set-TextValue "E47" '  -7.46%  '
Set-TextValue "D48" '0.0898'
Set-TextValue "E48" '  -4.15%  '
Set-TextValue "E49" '  -5.34%  '
Set-TextValue "D50" '16.78'
Set-TextValue "E50" '  -4.35%  '
Set-TextValue "E51" '  -5.21%  '
